$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for cells whose new content could otherwise be
# auto-coerced by Excel into a number/date (e.g. "1.00", "2.301.54").
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "55.061.42"
$ws.Range("E2").Value = "  +1.06%  "
Set-TextValue "D3" "2.301.54"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue "D5" "508.53"
$ws.Range("E5").Value = "  +0.82%  "
Set-TextValue "D6" "130.33"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  -0.42%  "
Set-TextValue "D8" "0.532"
$ws.Range("E8").Value = "  +0.49%  "
Set-TextValue "D9" "2.332.88"
$ws.Range("E9").Value = "  +1.54%  "
Set-TextValue "D10" "0.0982"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("E11").Value = "  +1.72%  "
Set-TextValue "D12" "5.07"
$ws.Range("E12").Value = "  +7.07%  "
$ws.Range("E13").Value = "  +1.57%  "
Set-TextValue "D14" "24.14"
$ws.Range("E14").Value = "  +4.69%  "
Set-TextValue "D15" "2.710.78"
$ws.Range("E15").Value = "  +0.64%  "
Set-TextValue "D16" "55.081.53"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("E17").Value = "  +1.49%  "
Set-TextValue "D18" "2.335.45"
$ws.Range("E18").Value = "  +0.76%  "
Set-TextValue "D19" "10.75"
$ws.Range("E19").Value = "  +4.05%  "
Set-TextValue "D20" "4.20"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "312.73"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D22" "6.69"
$ws.Range("E22").Value = "  +4.04%  "
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  +0.11%  "
Set-TextValue "D24" "60.76"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  -0.70%  "
Set-TextValue "D26" "0.152"
$ws.Range("E26").Value = "  +0.02%  "
Set-TextValue "D27" "7.55"
$ws.Range("E27").Value = "  +2.62%  "
Set-TextValue "D28" "173.10"
$ws.Range("E28").Value = "  -0.18%  "
Set-TextValue "D29" "6.19"
$ws.Range("E29").Value = "  +2.95%  "
Set-TextValue "D30" "0.0₃0713"
$ws.Range("E30").Value = "  +2.55%  "
Set-TextValue "D31" "1.16"
$ws.Range("E31").Value = "  +4.81%  "
$ws.Range("E32").Value = "  +0.42%  "
Set-TextValue "D33" "18.14"
$ws.Range("E33").Value = "  +1.60%  "
Set-TextValue "D34" "0.998"
$ws.Range("E34").Value = "  -0.06%  "
Set-TextValue "D35" "0.994"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D36" "1.24"
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D37" "0.922"
$ws.Range("E37").Value = "  -6.10%  "
Set-TextValue "D38" "3.91"
$ws.Range("E38").Value = "  +3.54%  "
Set-TextValue "D39" "36.81"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "1.45"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D41" "0.379"
$ws.Range("E41").Value = "  +1.08%  "
Set-TextValue "D42" "135.76"
$ws.Range("E42").Value = "  +8.26%  "
Set-TextValue "D43" "5.16"
$ws.Range("E43").Value = "  +5.10%  "
Set-TextValue "D44" "3.45"
$ws.Range("E44").Value = "  +1.39%  "
Set-TextValue "D45" "261.42"
$ws.Range("E45").Value = "  +6.39%  "
Set-TextValue "D46" "0.0505"
$ws.Range("E46").Value = "  +1.52%  "
Set-TextValue "D47" "0.0913"
$ws.Range("E47").Value = "  +2.01%  "
Set-TextValue "D48" "0.558"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("E49").Value = "  +1.15%  "
Set-TextValue "D50" "0.0211"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("E51").Value = "  +0.35%  "
